# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.231.42'
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").Value = '1.871.81'
$ws.Range("E3").Value = '  +3.37%  '

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = "'311.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.21%  '

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").Value = "'0.5015"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.26%  '

$ws.Range("D8").Value = "'0.3955"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.77%  '

$ws.Range("D9").Value = "'0.09871"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +26.15%  '

$ws.Range("D10").Value = "'1.138"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.61%  '

$ws.Range("D11").Value = "'41.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.43%  '

$ws.Range("D12").Value = "'6.474"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.49%  '

$ws.Range("D13").Value = "'21.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.76%  '

$ws.Range("D14").Value = '1.871.32'
$ws.Range("E14").Value = '  +3.56%  '

$ws.Range("D15").Value = "'1.002"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.02%  '

$ws.Range("D16").Value = "'7.402"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.76%  '

$ws.Range("D17").Value = "'0.00001137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.40%  '

$ws.Range("D18").Value = "'93.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.87%  '

$ws.Range("D19").Value = "'0.06643"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.91%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = "'17.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.75%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("D22").Value = "'6.117"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.70%  '

$ws.Range("D23").Value = '28.283.08'
$ws.Range("E23").Value = '  -0.13%  '

$ws.Range("D24").Value = "'11.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.88%  '

$ws.Range("D25").Value = "'2.274"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.76%  '

$ws.Range("D26").Value = "'2.560"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.89%  '

$ws.Range("D27").Value = "'21.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.82%  '

$ws.Range("D28").Value = '2.082.90'
$ws.Range("E28").Value = '  +3.29%  '

$ws.Range("D29").Value = "'158.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.12%  '

$ws.Range("D30").Value = "'128.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.31%  '

$ws.Range("D31").Value = "'0.1060"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.43%  '

$ws.Range("D32").Value = "'1.059"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.05%  '

$ws.Range("D33").Value = "'5.636"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.05%  '

$ws.Range("D34").Value = "'3.608"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.16%  '

$ws.Range("D35").Value = "'0.06809"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.15%  '

$ws.Range("D36").Value = "'9.481"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.57%  '

$ws.Range("D37").Value = "'0.02404"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.11%  '

$ws.Range("D38").Value = "'0.2187"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.39%  '

$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = "'11.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.98%  '

$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = "'5.016"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.11%  '

$ws.Range("D41").Value = "'0.6307"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.01%  '

$ws.Range("D42").Value = "'1.175"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.35%  '

$ws.Range("D44").Value = "'13.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.94%  '

$ws.Range("D45").Value = "'0.6019"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.89%  '

$ws.Range("D46").Value = "'3.665"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.72%  '

$ws.Range("D47").Value = "'1.272"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.55%  '

$ws.Range("D48").Value = "'124.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.37%  '

$ws.Range("D49").Value = "'1.992"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.62%  '

$ws.Range("D50").Value = "'1.202"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.76%  '

$ws.Range("D51").Value = "'1.125"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.64%  '

